$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L3").Value = "*maa://22880 (69.57), maa://20276 (83.45), *maa://22749 (66.67)"
$ws.Range("P3").Value = "maa://21249 (95.24), maa://26254 (95.65)"
$ws.Range("X3").Value = "maa://27396 (85.47), maa://27484 (95.79), maa://27480 (82.35)"
$ws.Range("T4").Value = "maa://32509 (97.73), maa://22754 (91.67), maa://27295 (81.82), *maa://21746 (55.81), *maa://31008 (78.05)"
$ws.Range("X4").Value = "**maa://32495 (47.01), ***maa://31785 (20.87), ***maa://36683 (28.26)"
$ws.Range("AF4").Value = "*maa://30062 (61.36), ***maa://26209 (13.04), *maa://39394 (75.0)"
$ws.Range("D5").Value = "maa://21245 (82.23), maa://22744 (83.33)"
$ws.Range("T6").Value = "maa://37411 (83.33)"
$ws.Range("X7").Value = "maa://22399 (94.78), *maa://22758 (70.91)"
$ws.Range("A8").Value = "更新日期：2024.11.03 01:13:15"
$ws.Range("L9").Value = "maa://22762 (91.57), maa://39552 (87.5)"
$ws.Range("P10").Value = "maa://28977 (93.59), *maa://23264 (61.82), maa://36669 (85.19)"
$ws.Range("X11").Value = "maa://36713 (97.93)"
$ws.Range("L14").Value = "maa://26245 (96.12), maa://21288 (96.21), maa://36682 (100.0), maa://39841 (93.88)"
$ws.Range("P14").Value = "maa://23250 (98.53), maa://20107 (87.1), maa://22772 (100.0), **maa://22745 (50.0)"
$ws.Range("T14").Value = "maa://22521 (94.57), maa://42751 (100.0)"
$ws.Range("T16").Value = "maa://22729 (95.24), *maa://28648 (69.64), maa://36674 (81.25)"
$ws.Range("H18").Value = "maa://24421 (90.41)"
$ws.Range("X18").Value = "maa://21917 (97.56), maa://22741 (83.33)"
$ws.Range("H21").Value = "maa://24372 (96.55)"
$ws.Range("X21").Value = "maa://20110 (86.76), maa://34946 (91.43)"
$ws.Range("L22").Value = "maa://27127 (85.71), *maa://22751 (76.19)"
$ws.Range("D24").Value = "maa://24368 (80.42)"
$ws.Range("H26").Value = "maa://24913 (91.78)"
$ws.Range("T28").Value = "maa://23263 (94.62), *maa://29765 (60.0)"
$ws.Range("P30").Value = "maa://21442 (99.5)"
$ws.Range("P33").Value = "*maa://21956 (79.1), maa://22730 (82.14)"
$ws.Range("P39").Value = "maa://24709 (92.38)"
$ws.Range("P41").Value = "**maa://35616 (38.24)"
$ws.Range("H44").Value = "maa://29768 (97.66), maa://27728 (96.0)"
$ws.Range("H59").Value = "maa://27746 (83.33), maa://31270 (95.19)"
